# Update the Tgfb2-Tgfbr2 LR-pairs sheet with the new TPM-based values.
# The data now also includes "Resolving-Mac" as a sending/target cluster
# (rows 14-17), and all numeric columns (G-T) were recomputed with the
# new TPM normalisation, extending the used range from A1:T13 to A1:T17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full replacement block for rows 2-17 (16 data rows x 20 cols,
# columns A-T) and write it in a single Range assignment.
$data = New-Object 'object[,]' 16,20
# Row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Tgfb2"
$data[0,2] = "Tgfbr2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 2.564153333333334
$data[0,7] = 7.692460000000001
$data[0,8] = 0.05249149542937438
$data[0,9] = 0.05249149542937438
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 27.85106533333333
$data[0,13] = 83.553196
$data[0,14] = 0.1861900221007236
$data[0,15] = 0.1861900221007236
$data[0,16] = 71.41440201135111
$data[0,17] = 642.7296181021601
$data[0,18] = 0.009773392694095249
$data[0,19] = 0.009773392694095249
# Row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Tgfb2"
$data[1,2] = "Tgfbr2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 2.564153333333334
$data[1,7] = 7.692460000000001
$data[1,8] = 0.05249149542937438
$data[1,9] = 0.05249149542937438
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 70.710031
$data[1,13] = 212.130093
$data[1,14] = 0.4727109026912454
$data[1,15] = 0.4727109026912454
$data[1,16] = 181.3113616887534
$data[1,17] = 1631.80225519878
$data[1,18] = 0.02481330218803295
$data[1,19] = 0.02481330218803295
# Row 4: ECs -> MuSCs
$data[2,0] = "ECs"
$data[2,1] = "Tgfb2"
$data[2,2] = "Tgfbr2"
$data[2,3] = "MuSCs"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 2.564153333333334
$data[2,7] = 7.692460000000001
$data[2,8] = 0.05249149542937438
$data[2,9] = 0.05249149542937438
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 15.018964
$data[2,13] = 45.056892
$data[2,14] = 0.1004048213460311
$data[2,15] = 0.1004048213460311
$data[2,16] = 38.51092660381333
$data[2,17] = 346.59833943432
$data[2,18] = 0.005270399220772341
$data[2,19] = 0.005270399220772341
# Row 5: ECs -> Resolving-Mac
$data[3,0] = "ECs"
$data[3,1] = "Tgfb2"
$data[3,2] = "Tgfbr2"
$data[3,3] = "Resolving-Mac"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.564153333333334
$data[3,7] = 7.692460000000001
$data[3,8] = 0.05249149542937438
$data[3,9] = 0.05249149542937438
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 36.00403133333333
$data[3,13] = 108.012094
$data[3,14] = 0.2406942538619999
$data[3,15] = 0.2406942538619999
$data[3,16] = 92.31985695680444
$data[3,17] = 830.87871261124
$data[3,18] = 0.01263440132647385
$data[3,19] = 0.01263440132647385
# Row 6: FAPs -> ECs
$data[4,0] = "FAPs"
$data[4,1] = "Tgfb2"
$data[4,2] = "Tgfbr2"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 16.59481266666667
$data[4,7] = 49.78443799999999
$data[4,8] = 0.3397170215679993
$data[4,9] = 0.3397170215679993
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 27.85106533333333
$data[4,13] = 83.553196
$data[4,14] = 0.1861900221007236
$data[4,15] = 0.1861900221007236
$data[4,16] = 462.1832117737609
$data[4,17] = 4159.648905963848
$data[4,18] = 0.0632519197537378
$data[4,19] = 0.06325191975373778
# Row 7: FAPs -> FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Tgfb2"
$data[5,2] = "Tgfbr2"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 16.59481266666667
$data[5,7] = 49.78443799999999
$data[5,8] = 0.3397170215679993
$data[5,9] = 0.3397170215679993
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 70.710031
$data[5,13] = 212.130093
$data[5,14] = 0.4727109026912454
$data[5,15] = 0.4727109026912454
$data[5,16] = 1173.419718099193
$data[5,17] = 10560.77746289273
$data[5,18] = 0.1605879399249903
$data[5,19] = 0.1605879399249902
# Row 8: FAPs -> MuSCs
$data[6,0] = "FAPs"
$data[6,1] = "Tgfb2"
$data[6,2] = "Tgfbr2"
$data[6,3] = "MuSCs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 16.59481266666667
$data[6,7] = 49.78443799999999
$data[6,8] = 0.3397170215679993
$data[6,9] = 0.3397170215679993
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 15.018964
$data[6,13] = 45.056892
$data[6,14] = 0.1004048213460311
$data[6,15] = 0.1004048213460311
$data[6,16] = 249.2368940274106
$data[6,17] = 2243.132046246696
$data[6,18] = 0.03410922685874075
$data[6,19] = 0.03410922685874074
# Row 9: FAPs -> Resolving-Mac
$data[7,0] = "FAPs"
$data[7,1] = "Tgfb2"
$data[7,2] = "Tgfbr2"
$data[7,3] = "Resolving-Mac"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 16.59481266666667
$data[7,7] = 49.78443799999999
$data[7,8] = 0.3397170215679993
$data[7,9] = 0.3397170215679993
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 36.00403133333333
$data[7,13] = 108.012094
$data[7,14] = 0.2406942538619999
$data[7,15] = 0.2406942538619999
$data[7,16] = 597.4801552214635
$data[7,17] = 5377.321396993171
$data[7,18] = 0.08176793503053052
$data[7,19] = 0.08176793503053051
# Row 10: MuSCs -> ECs
$data[8,0] = "MuSCs"
$data[8,1] = "Tgfb2"
$data[8,2] = "Tgfbr2"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 29.664466
$data[8,7] = 88.993398
$data[8,8] = 0.6072695268303631
$data[8,9] = 0.6072695268303631
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 27.85106533333333
$data[8,13] = 83.553196
$data[8,14] = 0.1861900221007236
$data[8,15] = 0.1861900221007236
$data[8,16] = 826.1869806444454
$data[8,17] = 7435.682825800008
$data[8,18] = 0.1130675266216413
$data[8,19] = 0.1130675266216413
# Row 11: MuSCs -> FAPs
$data[9,0] = "MuSCs"
$data[9,1] = "Tgfb2"
$data[9,2] = "Tgfbr2"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 29.664466
$data[9,7] = 88.993398
$data[9,8] = 0.6072695268303631
$data[9,9] = 0.6072695268303631
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 70.710031
$data[9,13] = 212.130093
$data[9,14] = 0.4727109026912454
$data[9,15] = 0.4727109026912454
$data[9,16] = 2097.575310458446
$data[9,17] = 18878.17779412601
$data[9,18] = 0.2870629262048664
$data[9,19] = 0.2870629262048664
# Row 12: MuSCs -> MuSCs
$data[10,0] = "MuSCs"
$data[10,1] = "Tgfb2"
$data[10,2] = "Tgfbr2"
$data[10,3] = "MuSCs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 29.664466
$data[10,7] = 88.993398
$data[10,8] = 0.6072695268303631
$data[10,9] = 0.6072695268303631
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 15.018964
$data[10,13] = 45.056892
$data[10,14] = 0.1004048213460311
$data[10,15] = 0.1004048213460311
$data[10,16] = 445.529546933224
$data[10,17] = 4009.765922399016
$data[10,18] = 0.06097278835029141
$data[10,19] = 0.06097278835029141
# Row 13: MuSCs -> Resolving-Mac
$data[11,0] = "MuSCs"
$data[11,1] = "Tgfb2"
$data[11,2] = "Tgfbr2"
$data[11,3] = "Resolving-Mac"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 29.664466
$data[11,7] = 88.993398
$data[11,8] = 0.6072695268303631
$data[11,9] = 0.6072695268303631
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 36.00403133333333
$data[11,13] = 108.012094
$data[11,14] = 0.2406942538619999
$data[11,15] = 0.2406942538619999
$data[11,16] = 1068.040363350601
$data[11,17] = 9612.363270155411
$data[11,18] = 0.146166285653564
$data[11,19] = 0.146166285653564
# Row 14: Resolving-Mac -> ECs
$data[12,0] = "Resolving-Mac"
$data[12,1] = "Tgfb2"
$data[12,2] = "Tgfbr2"
$data[12,3] = "ECs"
$data[12,4] = 1
$data[12,5] = 0.3333333333333333
$data[12,6] = 0.025497
$data[12,7] = 0.076491
$data[12,8] = 0.000521956172263265
$data[12,9] = 0.000521956172263265
$data[12,10] = 3
$data[12,11] = 1
$data[12,12] = 27.85106533333333
$data[12,13] = 83.553196
$data[12,14] = 0.1861900221007236
$data[12,15] = 0.1861900221007236
$data[12,16] = 0.7101186128040001
$data[12,17] = 6.391067515236
$data[12,18] = [double]"9.718303124930641E-05"
$data[12,19] = [double]"9.718303124930641E-05"
# Row 15: Resolving-Mac -> FAPs
$data[13,0] = "Resolving-Mac"
$data[13,1] = "Tgfb2"
$data[13,2] = "Tgfbr2"
$data[13,3] = "FAPs"
$data[13,4] = 1
$data[13,5] = 0.3333333333333333
$data[13,6] = 0.025497
$data[13,7] = 0.076491
$data[13,8] = 0.000521956172263265
$data[13,9] = 0.000521956172263265
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 70.710031
$data[13,13] = 212.130093
$data[13,14] = 0.4727109026912454
$data[13,15] = 0.4727109026912454
$data[13,16] = 1.802893660407
$data[13,17] = 16.226042943663
$data[13,18] = 0.0002467343733558352
$data[13,19] = 0.0002467343733558352
# Row 16: Resolving-Mac -> MuSCs
$data[14,0] = "Resolving-Mac"
$data[14,1] = "Tgfb2"
$data[14,2] = "Tgfbr2"
$data[14,3] = "MuSCs"
$data[14,4] = 1
$data[14,5] = 0.3333333333333333
$data[14,6] = 0.025497
$data[14,7] = 0.076491
$data[14,8] = 0.000521956172263265
$data[14,9] = 0.000521956172263265
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 15.018964
$data[14,13] = 45.056892
$data[14,14] = 0.1004048213460311
$data[14,15] = 0.1004048213460311
$data[14,16] = 0.382938525108
$data[14,17] = 3.446446725972
$data[14,18] = [double]"5.240691622655133E-05"
$data[14,19] = [double]"5.240691622655133E-05"
# Row 17: Resolving-Mac -> Resolving-Mac
$data[15,0] = "Resolving-Mac"
$data[15,1] = "Tgfb2"
$data[15,2] = "Tgfbr2"
$data[15,3] = "Resolving-Mac"
$data[15,4] = 1
$data[15,5] = 0.3333333333333333
$data[15,6] = 0.025497
$data[15,7] = 0.076491
$data[15,8] = 0.000521956172263265
$data[15,9] = 0.000521956172263265
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 36.00403133333333
$data[15,13] = 108.012094
$data[15,14] = 0.2406942538619999
$data[15,15] = 0.2406942538619999
$data[15,16] = 0.917994786906
$data[15,17] = 8.261953082153999
$data[15,18] = 0.0001256318514315721
$data[15,19] = 0.0001256318514315721

$ws.Range("A2:T17").Value = $data

